# Project To Project Transfer.xlsx - System Setup changes
# - Clear the "Transaction Date" values (6/29/2021) from rows 3 and 5,
#   which removes that now-unused string from the shared strings table.
# - Select row 7 downward (whole-row selection), matching the saved
#   sheet view selection of A7:XFD1048576.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the date cells that are no longer needed.
$ws.Range("B3").ClearContents()
$ws.Range("B5").ClearContents()

# Update the saved selection to match rows 7 through the bottom of the sheet.
$ws.Rows("7:1048576").Select()

$wb.Save()
